$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that held the "5840560 - Marco Antonio Carvalho Pereira" entry
# on its own (row 13); everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) now shows the teacher's name instead of the long
# "Desenvolver um projeto..." paragraph.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Row 13 (Programa resumido:) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date instead of the long
# "Noções de Gestão de Projetos..." syllabus text.
$ws.Range("B15").Value = "01/01/2015"
$ws.Range("C15").Value = "01/01/2015"

# Row 18 (Método:) now shows the teacher's name.
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Row 19 (Critério:) now holds the "método utilizado..." paragraph that used
# to belong to Método:.
$metodoTexto = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.
Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto

# Row 20 (Norma de recuperação:) now holds the grading-criteria paragraph
# that used to belong to Critério:.
$notaTexto = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("B20").Value = $notaTexto
$ws.Range("C20").Value = $notaTexto

# Row 21 (Bibliografia:) now just says "Não há recuperação".
$ws.Range("B21").Value = "Não há recuperação"
$ws.Range("C21").Value = "Não há recuperação"
